$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2663.875
$ws.Range("I19").Value = 3476.5
$ws.Range("K19").Value = 3476.5
$ws.Range("M19").Value = -3301.5
$ws.Range("H69").Value = 2033.3334
$ws.Range("J69").Value = 1500
$ws.Range("L69").Value = 4500
$ws.Range("N69").Value = -6248
$ws.Range("H72").Value = 2033.3334
$ws.Range("J72").Value = 1500
$ws.Range("L72").Value = 13500
$ws.Range("N72").Value = -22236
$ws.Range("H96").Value = 12500724
$ws.Range("I96").Value = 25000418
$ws.Range("J96").Value = 1029.9
$ws.Range("K96").Value = 75001254
$ws.Range("L96").Value = 3089.7
$ws.Range("M96").Value = -74999881
$ws.Range("N96").Value = -5835.700000000001
$ws.Range("H129").Value = 830.16
$ws.Range("I129").Value = 471
$ws.Range("J129").Value = 888.6279
$ws.Range("K129").Value = 1413
$ws.Range("L129").Value = 2665.8837
$ws.Range("M129").Value = 3587
$ws.Range("N129").Value = -12665.8837
$ws.Range("H135").Value = 26324680
$ws.Range("J135").Value = 83359000
$ws.Range("L135").Value = 750231000
$ws.Range("N135").Value = -750236070
$ws.Range("H137").Value = 2452.4358
$ws.Range("I137").Value = 2608.348
$ws.Range("J137").Value = 2228.3125
$ws.Range("K137").Value = 7825.044
$ws.Range("L137").Value = 6684.9375
$ws.Range("M137").Value = -5275.044
$ws.Range("N137").Value = -11784.9375
$ws.Range("H138").Value = 2341
$ws.Range("I138").Value = 666.94446
$ws.Range("J138").Value = 3417.1785
$ws.Range("K138").Value = 2000.83338
$ws.Range("L138").Value = 10251.5355
$ws.Range("M138").Value = 3139.16662
$ws.Range("N138").Value = -20531.5355

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1910.081
$ws.Range("I32").Value = 1578.5469
$ws.Range("K32").Value = 1578.5469
$ws.Range("M32").Value = -1291.5469
$ws.Range("H61").Value = 502581.88
$ws.Range("I61").Value = 751175.8
$ws.Range("K61").Value = 751175.8
$ws.Range("M61").Value = -750963.8
$ws.Range("H74").Value = 3635.2632
$ws.Range("I74").Value = 4750.909
$ws.Range("J74").Value = 2101.25
$ws.Range("K74").Value = 4750.909
$ws.Range("L74").Value = 2101.25
$ws.Range("M74").Value = -3876.909
$ws.Range("N74").Value = -3849.25
$ws.Range("H77").Value = 3635.2632
$ws.Range("I77").Value = 4750.909
$ws.Range("J77").Value = 2101.25
$ws.Range("K77").Value = 23754.545
$ws.Range("L77").Value = 10506.25
$ws.Range("M77").Value = -19386.545
$ws.Range("N77").Value = -19242.25
$ws.Range("H102").Value = 8399.799999999999
$ws.Range("I102").Value = 7333
$ws.Range("K102").Value = 7333
$ws.Range("M102").Value = -5711
$ws.Range("H110").Value = 2618.5
$ws.Range("I110").Value = 1667.9
$ws.Range("J110").Value = 4995
$ws.Range("K110").Value = 1667.9
$ws.Range("L110").Value = 4995
$ws.Range("M110").Value = 377.0999999999999
$ws.Range("N110").Value = -9085
$ws.Range("H115").Value = 38839.5
$ws.Range("J115").Value = 38839.5
$ws.Range("L115").Value = 38839.5
$ws.Range("N115").Value = -41973.5
$ws.Range("H132").Value = 22532.541
$ws.Range("I132").Value = 1494.3334
$ws.Range("J132").Value = 169800
$ws.Range("K132").Value = 4483.0002
$ws.Range("L132").Value = 509400
$ws.Range("M132").Value = -1953.0002
$ws.Range("N132").Value = -514460
$ws.Range("H136").Value = 502581.88
$ws.Range("I136").Value = 751175.8
$ws.Range("K136").Value = 2253527.4
$ws.Range("M136").Value = -2250977.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3356.08
$ws.Range("I94").Value = 1287.0769
$ws.Range("J94").Value = 5597.5
$ws.Range("K94").Value = 1287.0769
$ws.Range("L94").Value = 5597.5
$ws.Range("M94").Value = -836.0769
$ws.Range("N94").Value = -6499.5
$ws.Range("H105").Value = 3573300
$ws.Range("I105").Value = 2020
$ws.Range("J105").Value = 5557344.5
$ws.Range("K105").Value = 2020
$ws.Range("L105").Value = 5557344.5
$ws.Range("M105").Value = -273
$ws.Range("N105").Value = -5560838.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9399.682000000001
$ws.Range("I31").Value = 22716.066
$ws.Range("J31").Value = 2511.8965
$ws.Range("K31").Value = 22716.066
$ws.Range("L31").Value = 2511.8965
$ws.Range("M31").Value = -22421.066
$ws.Range("N31").Value = -3101.8965
$ws.Range("H34").Value = 9399.682000000001
$ws.Range("I34").Value = 22716.066
$ws.Range("J34").Value = 2511.8965
$ws.Range("K34").Value = 22716.066
$ws.Range("L34").Value = 2511.8965
$ws.Range("M34").Value = -22514.066
$ws.Range("N34").Value = -2915.8965
$ws.Range("H122").Value = 6000.5
$ws.Range("I122").Value = 6000.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18001.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15551.5
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 14524.45
$ws.Range("I132").Value = 16090.441
$ws.Range("J132").Value = 5650.5
$ws.Range("K132").Value = 48271.323
$ws.Range("L132").Value = 16951.5
$ws.Range("M132").Value = -45741.323
$ws.Range("N132").Value = -22011.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 687.87177
$ws.Range("J5").Value = 706.7692
$ws.Range("L5").Value = 2120.3076
$ws.Range("N5").Value = -2344.3076
$ws.Range("H68").Value = 5066.4614
$ws.Range("I68").Value = 829.8
$ws.Range("J68").Value = 6075.1904
$ws.Range("K68").Value = 2489.4
$ws.Range("L68").Value = 18225.5712
$ws.Range("M68").Value = -1678.4
$ws.Range("N68").Value = -19847.5712
$ws.Range("H71").Value = 5066.4614
$ws.Range("I71").Value = 829.8
$ws.Range("J71").Value = 6075.1904
$ws.Range("K71").Value = 7468.2
$ws.Range("L71").Value = 54676.7136
$ws.Range("M71").Value = -3412.2
$ws.Range("N71").Value = -62788.7136
$ws.Range("H113").Value = 713.1667
$ws.Range("I113").Value = 699.5
$ws.Range("K113").Value = 2098.5
$ws.Range("M113").Value = 71.5
$ws.Range("H131").Value = 827.3
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 837.42267
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2512.26801
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12592.26801
$ws.Range("H132").Value = 933.3333
$ws.Range("I132").Value = 1150
$ws.Range("J132").Value = 825
$ws.Range("K132").Value = 10350
$ws.Range("L132").Value = 7425
$ws.Range("M132").Value = -7820
$ws.Range("N132").Value = -12485
$ws.Range("H135").Value = 687.87177
$ws.Range("J135").Value = 706.7692
$ws.Range("L135").Value = 6360.922799999999
$ws.Range("N135").Value = -11430.9228

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 29533.525
$ws.Range("I132").Value = 3209.1333
$ws.Range("J132").Value = 128250
$ws.Range("K132").Value = 9627.3999
$ws.Range("L132").Value = 384750
$ws.Range("M132").Value = -7097.3999
$ws.Range("N132").Value = -389810
$ws.Range("H138").Value = 58429
$ws.Range("J138").Value = 58429
$ws.Range("L138").Value = 58429
$ws.Range("N138").Value = -68709

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5666.3335
$ws.Range("I16").Value = 5666.3335
$ws.Range("K16").Value = 5666.3335
$ws.Range("M16").Value = -5496.3335
$ws.Range("H46").Value = 810.8
$ws.Range("I46").Value = 589.38464
$ws.Range("K46").Value = 589.38464
$ws.Range("M46").Value = -401.38464
$ws.Range("H61").Value = 4749.7393
$ws.Range("I61").Value = 2540.3635
$ws.Range("J61").Value = 6775
$ws.Range("K61").Value = 2540.3635
$ws.Range("L61").Value = 6775
$ws.Range("M61").Value = -2338.3635
$ws.Range("N61").Value = -7179
$ws.Range("H93").Value = 1825.2273
$ws.Range("J93").Value = 1364.8
$ws.Range("L93").Value = 1364.8
$ws.Range("N93").Value = -3860.8
$ws.Range("H102").Value = 30000
$ws.Range("I102").Value = 30000
$ws.Range("K102").Value = 30000
$ws.Range("M102").Value = -26755
$ws.Range("H113").Value = 4749.7393
$ws.Range("I113").Value = 2540.3635
$ws.Range("J113").Value = 6775
$ws.Range("K113").Value = 2540.3635
$ws.Range("L113").Value = 6775
$ws.Range("M113").Value = -370.3634999999999
$ws.Range("N113").Value = -11115
$ws.Range("H122").Value = 2991.2222
$ws.Range("I122").Value = 2230.6365
$ws.Range("K122").Value = 6691.9095
$ws.Range("M122").Value = -4241.9095
$ws.Range("H132").Value = 1574.5416
$ws.Range("I132").Value = 1338.95
$ws.Range("J132").Value = 2752.5
$ws.Range("K132").Value = 4016.85
$ws.Range("L132").Value = 8257.5
$ws.Range("M132").Value = -1486.85
$ws.Range("N132").Value = -13317.5
$ws.Range("H136").Value = 2539.4
$ws.Range("I136").Value = 1888.7
$ws.Range("J136").Value = 3840.8
$ws.Range("K136").Value = 5666.1
$ws.Range("L136").Value = 11522.4
$ws.Range("M136").Value = -3116.1
$ws.Range("N136").Value = -16622.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 31139.834
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 31139.834
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 31139.834
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -31277.834
$ws.Range("H62").Value = 5423
$ws.Range("I62").Value = 5480
$ws.Range("J62").Value = 5404
$ws.Range("K62").Value = 5480
$ws.Range("L62").Value = 5404
$ws.Range("M62").Value = -4856
$ws.Range("N62").Value = -6652
$ws.Range("H65").Value = 5423
$ws.Range("I65").Value = 5480
$ws.Range("J65").Value = 5404
$ws.Range("K65").Value = 27400
$ws.Range("L65").Value = 27020
$ws.Range("M65").Value = -24280
$ws.Range("N65").Value = -33260
$ws.Range("H132").Value = 3150
$ws.Range("I132").Value = 3111.111
$ws.Range("K132").Value = 9333.332999999999
$ws.Range("M132").Value = -6803.332999999999
$ws.Range("H136").Value = 2007.8334
$ws.Range("I136").Value = 1398
$ws.Range("J136").Value = 2443.4285
$ws.Range("K136").Value = 4194
$ws.Range("L136").Value = 7330.2855
$ws.Range("M136").Value = -1644
$ws.Range("N136").Value = -12430.2855
